# ksom_wallFollow_01_hold_02.xlsx -- "adjusted xls of motor faiulre files"
#
# The sheet formerly called "acc_max" is renamed to "acc_best", the
# previously-active sheet (F1S) is left with a lingering cell selection at
# F16, and "acc_best" becomes the newly active/selected sheet & tab when the
# file is saved (mirroring the workbookView/activeTab + per-sheet
# tabSelected/selection bookkeeping Excel performs automatically).

$wb = $excel.ActiveWorkbook

# 1) Rename "acc_max" -> "acc_best" (sheetId stays the same, only the name
#    changes).
$accSheet = $wb.Worksheets.Item("acc_max")
$accSheet.Name = "acc_best"

# 2) Leave a cell selection behind on the F1S sheet at F16 (this is recorded
#    in that sheet's <selection> element) without leaving it as the active
#    tab.
$f1sSheet = $wb.Worksheets.Item("F1S")
$f1sSheet.Range("F16").Select() | Out-Null

# 3) Make "acc_best" (renamed acc_max) the active sheet/tab -- this updates
#    workbook.xml's bookViews/workbookView@activeTab as well as the
#    sheetView@tabSelected flag on the relevant sheets.
$wb.Worksheets.Item("acc_best").Activate() | Out-Null

# 4) Restore/record the window geometry (minimized, position & size) as it
#    was when the workbook was last saved.
$win = $wb.Windows.Item(1)
$win.WindowState = -4140
$win.Left = 20175
$win.Top = 1560
$win.Width = 11520
$win.Height = 7875
